$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Homework 4 (row 9) got graded
$ws.Range("E9").Value = 1.0

# Midterm 2 (row 21) got graded: total points 75, score flag set to 1
$ws.Range("D21").Value = 75.0
# Give D21 the same formatting as the other "Your Score" input cells (e.g. E21)
$ws.Range("E21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E21").Value = 1.0

# Midterm 2 raw scores entered in column J for each gradeline row (10-21)
$ws.Range("J10").Value = 60.0
$ws.Range("J11").Value = 53.0
$ws.Range("J12").Value = 46.0
$ws.Range("J13").Value = 39.0
$ws.Range("J14").Value = 32.0
$ws.Range("J15").Value = 25.0
$ws.Range("J16").Value = 20.0
$ws.Range("J17").Value = 15.0
$ws.Range("J18").Value = 10.0
$ws.Range("J19").Value = 7.0
$ws.Range("J20").Value = 4.0
$ws.Range("J21").Value = 0.0

# Force recalculation of all formula cells that depend on the above
# (re-entering each formula ensures dependents are refreshed)
$cell = $ws.Range("D2")
$cell.Formula = $cell.Formula

$cell = $ws.Range("D16")
$cell.Formula = $cell.Formula

foreach ($row in 10..21) {
    $cell = $ws.Range("H$row")
    $cell.Formula = $cell.Formula
}

foreach ($row in 10..22) {
    $cell = $ws.Range("M$row")
    $cell.Formula = $cell.Formula
}

$cell = $ws.Range("I23")
$cell.Formula = $cell.Formula
$cell = $ws.Range("J23")
$cell.Formula = $cell.Formula
